$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N2 value change: QADB -> DIGITAL_CHANNEL_SEC
$ws.Range("N2").Value = "DIGITAL_CHANNEL_SEC"

# C1 header rename: DB_Value -> status_query2
$ws.Range("C1").Value = "status_query2"

# C2 value change: QAT_BPS -> new BEGIN..END SQL block
$ws.Range("C2").Value = "BEGIN UPDATE DC_SCHEDULED_TRAN_MASTER STM SET STM.STATE = 46 , STM.IS_DELETED = 1 WHERE STM.BILL_BENEFICIARY_ID = (SELECT BPB.BENEFICIARY_ID FROM DC_BILL_PAYMENT_BENEFICIARY BPB WHERE BPB.CONSUMER_NUMBER = '0400000069505' AND BPB.CUSTOMER_INFO_ID = (SELECT CI.CUSTOMER_INFO_ID FROM DC_CUSTOMER_INFO CI WHERE CI.CUSTOMER_NAME = 'ABBY') AND BPB.IS_ACTIVE = 1);COMMIT;END;"

# Column C width adjustment (no longer best-fit, now matches column B's width)
$ws.Columns("C").ColumnWidth = 18

# Selection moves to C2
$null = $ws.Range("C2").Select()
